$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update container_cost value (B2) per the Demand Model integration
$ws.Range("B2").Value = 3.6413454393883802

# Recalculate so dependent formulas (e.g. Total_Cost in B8) refresh their cached values
$excel.Calculate()
